# Add a new "Driver" worksheet at the end of the workbook (after "Command")
# containing a small lookup table of driver function names, matching the
# "add driver fop by xing" commit.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last sheet so it lands at the end
# of the tab strip (Login, Command, Driver).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Driver"

# Header row: device/module names across columns B:F
$ws.Range("B1").Value = "LCD"
$ws.Range("C1").Value = "LED"
$ws.Range("D1").Value = "BUZZER"
$ws.Range("E1").Value = "BUTTON"
$ws.Range("F1").Value = "CAMERA"

# Row labels in column A
$ws.Range("A2").Value = "func_trunon"
$ws.Range("A3").Value = "func_trunoff"

# "turn on" function names
$ws.Range("B2").Value = "lcd_on"
$ws.Range("C2").Value = "led_on"
$ws.Range("D2").Value = "buz_on"
$ws.Range("E2").Value = "btn_on"
$ws.Range("F2").Value = "cam_on"

# "turn off" function names
$ws.Range("B3").Value = "lcd_off"
$ws.Range("C3").Value = "led_off"
$ws.Range("D3").Value = "buz_off"
$ws.Range("E3").Value = "btn_off"
$ws.Range("F3").Value = "cam_off"

# Match the author's final selection/active sheet state.
[void]$ws.Range("F21").Select()

Write-Host "Driver sheet added"
